$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 30-123 with shifted values (new row n = old row n-2 for most data columns,
# rows 30-31 get brand-new values). Only cells whose value actually changes are written.

# Row 30
$ws.Cells.Item(30, 4).Value = 44620
$ws.Cells.Item(30, 13).Value = 5
$ws.Cells.Item(30, 14).Value = 360000
$ws.Cells.Item(30, 15).Value = 360000
$ws.Cells.Item(30, 16).Value = 360000
$ws.Cells.Item(30, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(30, 19).Value = 800
$ws.Cells.Item(30, 20).Value = 450

# Row 31
$ws.Cells.Item(31, 4).Value = 44620
$ws.Cells.Item(31, 13).Value = 65
$ws.Cells.Item(31, 14).Value = 16000
$ws.Cells.Item(31, 15).Value = 16000
$ws.Cells.Item(31, 16).Value = 16000
$ws.Cells.Item(31, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(31, 19).Value = 889

# Row 32
$ws.Cells.Item(32, 4).Value = 44308
$ws.Cells.Item(32, 13).Value = 150
$ws.Cells.Item(32, 15).Value = 15000
$ws.Cells.Item(32, 16).Value = 14067
$ws.Cells.Item(32, 19).Value = 782

# Row 33
$ws.Cells.Item(33, 4).Value = 44308
$ws.Cells.Item(33, 12).Value = "Primera"
$ws.Cells.Item(33, 13).Value = 90
$ws.Cells.Item(33, 14).Value = 17000
$ws.Cells.Item(33, 15).Value = 17000
$ws.Cells.Item(33, 16).Value = 17000
$ws.Cells.Item(33, 17).Value = "`$/caja 18 kilos empedrada"
$ws.Cells.Item(33, 19).Value = 944

# Row 34
$ws.Cells.Item(34, 4).Value = 44427
$ws.Cells.Item(34, 13).Value = 80
$ws.Cells.Item(34, 14).Value = 13000
$ws.Cells.Item(34, 16).Value = 13000
$ws.Cells.Item(34, 19).Value = 722

# Row 35
$ws.Cells.Item(35, 4).Value = 44453
$ws.Cells.Item(35, 12).Value = "Especial"
$ws.Cells.Item(35, 13).Value = 55
$ws.Cells.Item(35, 14).Value = 20000
$ws.Cells.Item(35, 15).Value = 20000
$ws.Cells.Item(35, 16).Value = 20000
$ws.Cells.Item(35, 19).Value = 1111

# Row 36
$ws.Cells.Item(36, 4).Value = 44340
$ws.Cells.Item(36, 13).Value = 105
$ws.Cells.Item(36, 14).Value = 12000
$ws.Cells.Item(36, 15).Value = 13000
$ws.Cells.Item(36, 16).Value = 12619
$ws.Cells.Item(36, 19).Value = 701

# Row 37
$ws.Cells.Item(37, 4).Value = 44434
$ws.Cells.Item(37, 13).Value = 40
$ws.Cells.Item(37, 14).Value = 13000
$ws.Cells.Item(37, 15).Value = 13000
$ws.Cells.Item(37, 16).Value = 13000
$ws.Cells.Item(37, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(37, 19).Value = 722
$ws.Cells.Item(37, 20).Value = 18

# Row 38
$ws.Cells.Item(38, 4).Value = 44270
$ws.Cells.Item(38, 13).Value = 65
$ws.Cells.Item(38, 14).Value = 15000
$ws.Cells.Item(38, 15).Value = 15000
$ws.Cells.Item(38, 16).Value = 15000
$ws.Cells.Item(38, 19).Value = 833

# Row 39
$ws.Cells.Item(39, 4).Value = 44270
$ws.Cells.Item(39, 13).Value = 5
$ws.Cells.Item(39, 14).Value = 250000
$ws.Cells.Item(39, 15).Value = 250000
$ws.Cells.Item(39, 16).Value = 250000
$ws.Cells.Item(39, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(39, 19).Value = 556
$ws.Cells.Item(39, 20).Value = 450

# Row 40
$ws.Cells.Item(40, 4).Value = 44391
$ws.Cells.Item(40, 14).Value = 12000
$ws.Cells.Item(40, 16).Value = 12500
$ws.Cells.Item(40, 19).Value = 694

# Row 41
$ws.Cells.Item(41, 4).Value = 44376
$ws.Cells.Item(41, 13).Value = 95

# Row 42
$ws.Cells.Item(42, 4).Value = 44307
$ws.Cells.Item(42, 14).Value = 13000
$ws.Cells.Item(42, 16).Value = 13000
$ws.Cells.Item(42, 19).Value = 722

# Row 43
$ws.Cells.Item(43, 4).Value = 44405
$ws.Cells.Item(43, 13).Value = 50

# Row 44
$ws.Cells.Item(44, 4).Value = 44377
$ws.Cells.Item(44, 12).Value = "Primera"
$ws.Cells.Item(44, 13).Value = 80
$ws.Cells.Item(44, 14).Value = 12000
$ws.Cells.Item(44, 15).Value = 13000
$ws.Cells.Item(44, 16).Value = 12500
$ws.Cells.Item(44, 19).Value = 694

# Row 45
$ws.Cells.Item(45, 4).Value = 44424
$ws.Cells.Item(45, 13).Value = 45
$ws.Cells.Item(45, 14).Value = 13000
$ws.Cells.Item(45, 15).Value = 13000
$ws.Cells.Item(45, 16).Value = 13000
$ws.Cells.Item(45, 19).Value = 722

# Row 46
$ws.Cells.Item(46, 4).Value = 44306
$ws.Cells.Item(46, 12).Value = "Especial"
$ws.Cells.Item(46, 13).Value = 40
$ws.Cells.Item(46, 14).Value = 18000
$ws.Cells.Item(46, 15).Value = 18000
$ws.Cells.Item(46, 16).Value = 18000
$ws.Cells.Item(46, 19).Value = 1000

# Row 47
$ws.Cells.Item(47, 4).Value = 44384
$ws.Cells.Item(47, 13).Value = 95
$ws.Cells.Item(47, 14).Value = 12000
$ws.Cells.Item(47, 15).Value = 12000
$ws.Cells.Item(47, 16).Value = 12000
$ws.Cells.Item(47, 19).Value = 667

# Row 48
$ws.Cells.Item(48, 4).Value = 44383
$ws.Cells.Item(48, 13).Value = 120
$ws.Cells.Item(48, 16).Value = 12542
$ws.Cells.Item(48, 19).Value = 697

# Row 49
$ws.Cells.Item(49, 4).Value = 44295
$ws.Cells.Item(49, 12).Value = "Primera"
$ws.Cells.Item(49, 13).Value = 40
$ws.Cells.Item(49, 14).Value = 13000
$ws.Cells.Item(49, 15).Value = 13000
$ws.Cells.Item(49, 16).Value = 13000
$ws.Cells.Item(49, 19).Value = 722

# Row 50
$ws.Cells.Item(50, 4).Value = 44369
$ws.Cells.Item(50, 13).Value = 105
$ws.Cells.Item(50, 14).Value = 12000
$ws.Cells.Item(50, 15).Value = 13000
$ws.Cells.Item(50, 16).Value = 12381
$ws.Cells.Item(50, 19).Value = 688

# Row 51
$ws.Cells.Item(51, 4).Value = 44369
$ws.Cells.Item(51, 12).Value = "Segunda"
$ws.Cells.Item(51, 13).Value = 65
$ws.Cells.Item(51, 14).Value = 8000
$ws.Cells.Item(51, 15).Value = 8000
$ws.Cells.Item(51, 16).Value = 8000
$ws.Cells.Item(51, 19).Value = 444

# Row 52
$ws.Cells.Item(52, 4).Value = 44278
$ws.Cells.Item(52, 13).Value = 80
$ws.Cells.Item(52, 15).Value = 15000
$ws.Cells.Item(52, 16).Value = 14500
$ws.Cells.Item(52, 19).Value = 806

# Row 53
$ws.Cells.Item(53, 4).Value = 44435
$ws.Cells.Item(53, 12).Value = "Primera"
$ws.Cells.Item(53, 13).Value = 140
$ws.Cells.Item(53, 14).Value = 13000
$ws.Cells.Item(53, 15).Value = 13000
$ws.Cells.Item(53, 16).Value = 13000
$ws.Cells.Item(53, 19).Value = 722

# Row 54
$ws.Cells.Item(54, 4).Value = 44292
$ws.Cells.Item(54, 13).Value = 50

# Row 55
$ws.Cells.Item(55, 4).Value = 44314
$ws.Cells.Item(55, 12).Value = "Especial"
$ws.Cells.Item(55, 13).Value = 35
$ws.Cells.Item(55, 14).Value = 20000
$ws.Cells.Item(55, 15).Value = 20000
$ws.Cells.Item(55, 16).Value = 20000
$ws.Cells.Item(55, 19).Value = 1111

# Row 56
$ws.Cells.Item(56, 4).Value = 44314
$ws.Cells.Item(56, 13).Value = 55
$ws.Cells.Item(56, 14).Value = 14000
$ws.Cells.Item(56, 15).Value = 14000
$ws.Cells.Item(56, 16).Value = 14000
$ws.Cells.Item(56, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(56, 19).Value = 778
$ws.Cells.Item(56, 20).Value = 18

# Row 57
$ws.Cells.Item(57, 4).Value = 44354
$ws.Cells.Item(57, 13).Value = 120
$ws.Cells.Item(57, 14).Value = 13000
$ws.Cells.Item(57, 15).Value = 14000
$ws.Cells.Item(57, 16).Value = 13542
$ws.Cells.Item(57, 19).Value = 752

# Row 58
$ws.Cells.Item(58, 4).Value = 44354
$ws.Cells.Item(58, 13).Value = 5
$ws.Cells.Item(58, 14).Value = 270000
$ws.Cells.Item(58, 15).Value = 270000
$ws.Cells.Item(58, 16).Value = 270000
$ws.Cells.Item(58, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(58, 19).Value = 600
$ws.Cells.Item(58, 20).Value = 450

# Row 59
$ws.Cells.Item(59, 4).Value = 44406
$ws.Cells.Item(59, 12).Value = "Primera"
$ws.Cells.Item(59, 13).Value = 70
$ws.Cells.Item(59, 14).Value = 12000
$ws.Cells.Item(59, 15).Value = 13000
$ws.Cells.Item(59, 16).Value = 12429
$ws.Cells.Item(59, 19).Value = 690

# Row 60
$ws.Cells.Item(60, 4).Value = 44260
$ws.Cells.Item(60, 13).Value = 65
$ws.Cells.Item(60, 14).Value = 15000
$ws.Cells.Item(60, 15).Value = 15000
$ws.Cells.Item(60, 16).Value = 15000
$ws.Cells.Item(60, 19).Value = 833

# Row 61
$ws.Cells.Item(61, 4).Value = 44341
$ws.Cells.Item(61, 12).Value = "Especial"
$ws.Cells.Item(61, 14).Value = 15000
$ws.Cells.Item(61, 15).Value = 15000
$ws.Cells.Item(61, 16).Value = 15000
$ws.Cells.Item(61, 19).Value = 833

# Row 62
$ws.Cells.Item(62, 4).Value = 44341
$ws.Cells.Item(62, 14).Value = 12000
$ws.Cells.Item(62, 15).Value = 12000
$ws.Cells.Item(62, 16).Value = 12000
$ws.Cells.Item(62, 19).Value = 667

# Row 63
$ws.Cells.Item(63, 4).Value = 44245
$ws.Cells.Item(63, 13).Value = 65
$ws.Cells.Item(63, 14).Value = 17000
$ws.Cells.Item(63, 15).Value = 17000
$ws.Cells.Item(63, 16).Value = 17000
$ws.Cells.Item(63, 19).Value = 944

# Row 64
$ws.Cells.Item(64, 4).Value = 44305
$ws.Cells.Item(64, 13).Value = 80
$ws.Cells.Item(64, 14).Value = 13000
$ws.Cells.Item(64, 15).Value = 13000
$ws.Cells.Item(64, 16).Value = 13000
$ws.Cells.Item(64, 19).Value = 722

# Row 65
$ws.Cells.Item(65, 4).Value = 44442
$ws.Cells.Item(65, 13).Value = 75
$ws.Cells.Item(65, 14).Value = 12000
$ws.Cells.Item(65, 15).Value = 12000
$ws.Cells.Item(65, 16).Value = 12000
$ws.Cells.Item(65, 19).Value = 667

# Row 66
$ws.Cells.Item(66, 4).Value = 44246
$ws.Cells.Item(66, 14).Value = 17000
$ws.Cells.Item(66, 15).Value = 17000
$ws.Cells.Item(66, 16).Value = 17000
$ws.Cells.Item(66, 19).Value = 944

# Row 67
$ws.Cells.Item(67, 4).Value = 44323
$ws.Cells.Item(67, 13).Value = 40
$ws.Cells.Item(67, 14).Value = 13000
$ws.Cells.Item(67, 16).Value = 13000
$ws.Cells.Item(67, 19).Value = 722

# Row 68
$ws.Cells.Item(68, 4).Value = 44398
$ws.Cells.Item(68, 13).Value = 55

# Row 69
$ws.Cells.Item(69, 4).Value = 44392
$ws.Cells.Item(69, 13).Value = 90
$ws.Cells.Item(69, 14).Value = 12000
$ws.Cells.Item(69, 16).Value = 12444
$ws.Cells.Item(69, 19).Value = 691

# Row 70
$ws.Cells.Item(70, 4).Value = 44328
$ws.Cells.Item(70, 13).Value = 65
$ws.Cells.Item(70, 14).Value = 13000
$ws.Cells.Item(70, 16).Value = 13000
$ws.Cells.Item(70, 19).Value = 722

# Row 71
$ws.Cells.Item(71, 4).Value = 44433
$ws.Cells.Item(71, 13).Value = 80
$ws.Cells.Item(71, 14).Value = 13000
$ws.Cells.Item(71, 15).Value = 13000
$ws.Cells.Item(71, 16).Value = 13000
$ws.Cells.Item(71, 19).Value = 722

# Row 72
$ws.Cells.Item(72, 4).Value = 44382
$ws.Cells.Item(72, 13).Value = 115
$ws.Cells.Item(72, 14).Value = 11000
$ws.Cells.Item(72, 16).Value = 12565
$ws.Cells.Item(72, 19).Value = 698

# Row 73
$ws.Cells.Item(73, 4).Value = 44265
$ws.Cells.Item(73, 13).Value = 40
$ws.Cells.Item(73, 14).Value = 15000
$ws.Cells.Item(73, 15).Value = 15000
$ws.Cells.Item(73, 16).Value = 15000
$ws.Cells.Item(73, 19).Value = 833

# Row 74
$ws.Cells.Item(74, 4).Value = 44363
$ws.Cells.Item(74, 13).Value = 80

# Row 75
$ws.Cells.Item(75, 4).Value = 44386
$ws.Cells.Item(75, 13).Value = 95

# Row 76
$ws.Cells.Item(76, 4).Value = 44322
$ws.Cells.Item(76, 12).Value = "Primera"
$ws.Cells.Item(76, 13).Value = 200
$ws.Cells.Item(76, 14).Value = 13000
$ws.Cells.Item(76, 15).Value = 13000
$ws.Cells.Item(76, 16).Value = 13000
$ws.Cells.Item(76, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(76, 19).Value = 722

# Row 77
$ws.Cells.Item(77, 4).Value = 44358
$ws.Cells.Item(77, 13).Value = 80
$ws.Cells.Item(77, 14).Value = 12000
$ws.Cells.Item(77, 15).Value = 12000
$ws.Cells.Item(77, 16).Value = 12000
$ws.Cells.Item(77, 19).Value = 667

# Row 78
$ws.Cells.Item(78, 4).Value = 44313
$ws.Cells.Item(78, 12).Value = "Especial"
$ws.Cells.Item(78, 13).Value = 125
$ws.Cells.Item(78, 14).Value = 20000
$ws.Cells.Item(78, 15).Value = 20000
$ws.Cells.Item(78, 16).Value = 20000
$ws.Cells.Item(78, 17).Value = "`$/caja 18 kilos empedrada"
$ws.Cells.Item(78, 19).Value = 1111

# Row 79
$ws.Cells.Item(79, 4).Value = 44244
$ws.Cells.Item(79, 13).Value = 45
$ws.Cells.Item(79, 14).Value = 17000
$ws.Cells.Item(79, 15).Value = 17000
$ws.Cells.Item(79, 16).Value = 17000
$ws.Cells.Item(79, 19).Value = 944

# Row 80
$ws.Cells.Item(80, 4).Value = 44356
$ws.Cells.Item(80, 13).Value = 110
$ws.Cells.Item(80, 14).Value = 12000
$ws.Cells.Item(80, 15).Value = 12000
$ws.Cells.Item(80, 16).Value = 12000
$ws.Cells.Item(80, 19).Value = 667

# Row 81
$ws.Cells.Item(81, 4).Value = 44302
$ws.Cells.Item(81, 13).Value = 95
$ws.Cells.Item(81, 14).Value = 11000
$ws.Cells.Item(81, 15).Value = 12000
$ws.Cells.Item(81, 16).Value = 11474
$ws.Cells.Item(81, 19).Value = 637

# Row 82
$ws.Cells.Item(82, 4).Value = 44291

# Row 83
$ws.Cells.Item(83, 4).Value = 44326
$ws.Cells.Item(83, 13).Value = 115
$ws.Cells.Item(83, 16).Value = 13565
$ws.Cells.Item(83, 19).Value = 754

# Row 84
$ws.Cells.Item(84, 4).Value = 44348
$ws.Cells.Item(84, 12).Value = "Primera"
$ws.Cells.Item(84, 13).Value = 40
$ws.Cells.Item(84, 14).Value = 13000
$ws.Cells.Item(84, 15).Value = 13000
$ws.Cells.Item(84, 16).Value = 13000
$ws.Cells.Item(84, 19).Value = 722

# Row 85
$ws.Cells.Item(85, 4).Value = 44281
$ws.Cells.Item(85, 13).Value = 60
$ws.Cells.Item(85, 14).Value = 13000
$ws.Cells.Item(85, 15).Value = 14000
$ws.Cells.Item(85, 16).Value = 13417
$ws.Cells.Item(85, 19).Value = 745

# Row 86
$ws.Cells.Item(86, 4).Value = 44271
$ws.Cells.Item(86, 12).Value = "Especial"
$ws.Cells.Item(86, 13).Value = 35
$ws.Cells.Item(86, 14).Value = 20000
$ws.Cells.Item(86, 15).Value = 20000
$ws.Cells.Item(86, 16).Value = 20000
$ws.Cells.Item(86, 19).Value = 1111

# Row 87
$ws.Cells.Item(87, 4).Value = 44271
$ws.Cells.Item(87, 13).Value = 95
$ws.Cells.Item(87, 14).Value = 15000
$ws.Cells.Item(87, 15).Value = 15000
$ws.Cells.Item(87, 16).Value = 15000
$ws.Cells.Item(87, 19).Value = 833

# Row 88
$ws.Cells.Item(88, 4).Value = 44420
$ws.Cells.Item(88, 13).Value = 65
$ws.Cells.Item(88, 15).Value = 13000
$ws.Cells.Item(88, 16).Value = 13000
$ws.Cells.Item(88, 19).Value = 722

# Row 89
$ws.Cells.Item(89, 4).Value = 44343
$ws.Cells.Item(89, 12).Value = "Primera"
$ws.Cells.Item(89, 13).Value = 75
$ws.Cells.Item(89, 14).Value = 13000
$ws.Cells.Item(89, 15).Value = 13000
$ws.Cells.Item(89, 16).Value = 13000
$ws.Cells.Item(89, 19).Value = 722

# Row 90
$ws.Cells.Item(90, 4).Value = 44315
$ws.Cells.Item(90, 13).Value = 85
$ws.Cells.Item(90, 14).Value = 13000
$ws.Cells.Item(90, 15).Value = 14000
$ws.Cells.Item(90, 16).Value = 13529
$ws.Cells.Item(90, 19).Value = 752

# Row 91
$ws.Cells.Item(91, 4).Value = 44315
$ws.Cells.Item(91, 12).Value = "Segunda"
$ws.Cells.Item(91, 13).Value = 55
$ws.Cells.Item(91, 14).Value = 9000
$ws.Cells.Item(91, 15).Value = 9000
$ws.Cells.Item(91, 16).Value = 9000
$ws.Cells.Item(91, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(91, 19).Value = 500
$ws.Cells.Item(91, 20).Value = 18

# Row 92
$ws.Cells.Item(92, 4).Value = 44336
$ws.Cells.Item(92, 13).Value = 80
$ws.Cells.Item(92, 15).Value = 12000
$ws.Cells.Item(92, 16).Value = 12000
$ws.Cells.Item(92, 19).Value = 667

# Row 93
$ws.Cells.Item(93, 4).Value = 44336
$ws.Cells.Item(93, 13).Value = 3
$ws.Cells.Item(93, 14).Value = 230000
$ws.Cells.Item(93, 15).Value = 230000
$ws.Cells.Item(93, 16).Value = 230000
$ws.Cells.Item(93, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(93, 19).Value = 511
$ws.Cells.Item(93, 20).Value = 450

# Row 94
$ws.Cells.Item(94, 4).Value = 44400
$ws.Cells.Item(94, 13).Value = 120
$ws.Cells.Item(94, 16).Value = 12542
$ws.Cells.Item(94, 19).Value = 697

# Row 95
$ws.Cells.Item(95, 4).Value = 44334
$ws.Cells.Item(95, 13).Value = 90

# Row 96
$ws.Cells.Item(96, 4).Value = 44319
$ws.Cells.Item(96, 13).Value = 180
$ws.Cells.Item(96, 14).Value = 12000
$ws.Cells.Item(96, 15).Value = 13000
$ws.Cells.Item(96, 16).Value = 12444
$ws.Cells.Item(96, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(96, 19).Value = 691
$ws.Cells.Item(96, 20).Value = 18

# Row 97
$ws.Cells.Item(97, 4).Value = 44280
$ws.Cells.Item(97, 13).Value = 200
$ws.Cells.Item(97, 14).Value = 12000
$ws.Cells.Item(97, 15).Value = 12000
$ws.Cells.Item(97, 16).Value = 12000
$ws.Cells.Item(97, 19).Value = 667

# Row 98
$ws.Cells.Item(98, 4).Value = 44280
$ws.Cells.Item(98, 13).Value = 3
$ws.Cells.Item(98, 14).Value = 240000
$ws.Cells.Item(98, 15).Value = 240000
$ws.Cells.Item(98, 16).Value = 240000
$ws.Cells.Item(98, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(98, 19).Value = 533
$ws.Cells.Item(98, 20).Value = 450

# Row 99
$ws.Cells.Item(99, 4).Value = 44362
$ws.Cells.Item(99, 13).Value = 40
$ws.Cells.Item(99, 14).Value = 13000
$ws.Cells.Item(99, 16).Value = 13000
$ws.Cells.Item(99, 19).Value = 722

# Row 100
$ws.Cells.Item(100, 4).Value = 44431
$ws.Cells.Item(100, 13).Value = 20
$ws.Cells.Item(100, 14).Value = 13000
$ws.Cells.Item(100, 15).Value = 13000
$ws.Cells.Item(100, 16).Value = 13000
$ws.Cells.Item(100, 19).Value = 722

# Row 101
$ws.Cells.Item(101, 4).Value = 44365
$ws.Cells.Item(101, 13).Value = 70
$ws.Cells.Item(101, 14).Value = 12000
$ws.Cells.Item(101, 16).Value = 12571
$ws.Cells.Item(101, 19).Value = 698

# Row 102
$ws.Cells.Item(102, 4).Value = 44357
$ws.Cells.Item(102, 12).Value = "Primera"
$ws.Cells.Item(102, 13).Value = 125
$ws.Cells.Item(102, 14).Value = 12000
$ws.Cells.Item(102, 15).Value = 12000
$ws.Cells.Item(102, 16).Value = 12000
$ws.Cells.Item(102, 19).Value = 667

# Row 103
$ws.Cells.Item(103, 4).Value = 44397
$ws.Cells.Item(103, 13).Value = 85
$ws.Cells.Item(103, 15).Value = 13000
$ws.Cells.Item(103, 16).Value = 13000
$ws.Cells.Item(103, 19).Value = 722

# Row 104
$ws.Cells.Item(104, 4).Value = 44446
$ws.Cells.Item(104, 12).Value = "Especial"
$ws.Cells.Item(104, 13).Value = 40
$ws.Cells.Item(104, 14).Value = 20000
$ws.Cells.Item(104, 15).Value = 20000
$ws.Cells.Item(104, 16).Value = 20000
$ws.Cells.Item(104, 19).Value = 1111

# Row 105
$ws.Cells.Item(105, 4).Value = 44329
$ws.Cells.Item(105, 13).Value = 115
$ws.Cells.Item(105, 15).Value = 14000
$ws.Cells.Item(105, 16).Value = 13435
$ws.Cells.Item(105, 19).Value = 746

# Row 106
$ws.Cells.Item(106, 4).Value = 44355
$ws.Cells.Item(106, 13).Value = 115
$ws.Cells.Item(106, 14).Value = 13000
$ws.Cells.Item(106, 15).Value = 14000
$ws.Cells.Item(106, 16).Value = 13435
$ws.Cells.Item(106, 19).Value = 746

# Row 107
$ws.Cells.Item(107, 4).Value = 44294
$ws.Cells.Item(107, 13).Value = 100
$ws.Cells.Item(107, 14).Value = 13000
$ws.Cells.Item(107, 15).Value = 13000
$ws.Cells.Item(107, 16).Value = 13000
$ws.Cells.Item(107, 19).Value = 722

# Row 108
$ws.Cells.Item(108, 4).Value = 44617
$ws.Cells.Item(108, 13).Value = 20
$ws.Cells.Item(108, 14).Value = 16000
$ws.Cells.Item(108, 15).Value = 16000
$ws.Cells.Item(108, 16).Value = 16000
$ws.Cells.Item(108, 19).Value = 889

# Row 109
$ws.Cells.Item(109, 4).Value = 44264
$ws.Cells.Item(109, 13).Value = 30
$ws.Cells.Item(109, 14).Value = 16000
$ws.Cells.Item(109, 15).Value = 16000
$ws.Cells.Item(109, 16).Value = 16000
$ws.Cells.Item(109, 19).Value = 889

# Row 110
$ws.Cells.Item(110, 4).Value = 44396
$ws.Cells.Item(110, 13).Value = 175
$ws.Cells.Item(110, 14).Value = 12000
$ws.Cells.Item(110, 15).Value = 13000
$ws.Cells.Item(110, 16).Value = 12457
$ws.Cells.Item(110, 19).Value = 692

# Row 111
$ws.Cells.Item(111, 4).Value = 44279
$ws.Cells.Item(111, 13).Value = 120
$ws.Cells.Item(111, 14).Value = 12000
$ws.Cells.Item(111, 16).Value = 12417
$ws.Cells.Item(111, 19).Value = 690

# Row 112
$ws.Cells.Item(112, 4).Value = 44330
$ws.Cells.Item(112, 13).Value = 95
$ws.Cells.Item(112, 14).Value = 13000
$ws.Cells.Item(112, 15).Value = 14000
$ws.Cells.Item(112, 16).Value = 13526
$ws.Cells.Item(112, 19).Value = 751

# Row 113
$ws.Cells.Item(113, 4).Value = 44301
$ws.Cells.Item(113, 13).Value = 45
$ws.Cells.Item(113, 14).Value = 13000
$ws.Cells.Item(113, 16).Value = 13000
$ws.Cells.Item(113, 19).Value = 722

# Row 114
$ws.Cells.Item(114, 4).Value = 44370
$ws.Cells.Item(114, 13).Value = 65
$ws.Cells.Item(114, 15).Value = 12000
$ws.Cells.Item(114, 16).Value = 12000
$ws.Cells.Item(114, 19).Value = 667

# Row 115
$ws.Cells.Item(115, 4).Value = 44385
$ws.Cells.Item(115, 13).Value = 285
$ws.Cells.Item(115, 14).Value = 12000
$ws.Cells.Item(115, 15).Value = 13000
$ws.Cells.Item(115, 16).Value = 12561
$ws.Cells.Item(115, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(115, 19).Value = 698

# Row 116
$ws.Cells.Item(116, 4).Value = 44413
$ws.Cells.Item(116, 15).Value = 13000
$ws.Cells.Item(116, 16).Value = 12579
$ws.Cells.Item(116, 19).Value = 699

# Row 117
$ws.Cells.Item(117, 4).Value = 44312
$ws.Cells.Item(117, 13).Value = 210
$ws.Cells.Item(117, 14).Value = 13000
$ws.Cells.Item(117, 15).Value = 14000
$ws.Cells.Item(117, 16).Value = 13262
$ws.Cells.Item(117, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(117, 19).Value = 737

# Row 118
$ws.Cells.Item(118, 4).Value = 44399
$ws.Cells.Item(118, 13).Value = 95
$ws.Cells.Item(118, 14).Value = 12000
$ws.Cells.Item(118, 15).Value = 12000
$ws.Cells.Item(118, 16).Value = 12000
$ws.Cells.Item(118, 19).Value = 667

# Row 119
$ws.Cells.Item(119, 4).Value = 44615
$ws.Cells.Item(119, 13).Value = 50
$ws.Cells.Item(119, 14).Value = 16000
$ws.Cells.Item(119, 15).Value = 16000
$ws.Cells.Item(119, 16).Value = 16000
$ws.Cells.Item(119, 19).Value = 889

# Row 120
$ws.Cells.Item(120, 4).Value = 44277
$ws.Cells.Item(120, 13).Value = 160
$ws.Cells.Item(120, 14).Value = 13000
$ws.Cells.Item(120, 15).Value = 15000
$ws.Cells.Item(120, 16).Value = 14000
$ws.Cells.Item(120, 19).Value = 778

# Row 121
$ws.Cells.Item(121, 4).Value = 44258
$ws.Cells.Item(121, 13).Value = 110
$ws.Cells.Item(121, 14).Value = 15000
$ws.Cells.Item(121, 15).Value = 15000
$ws.Cells.Item(121, 16).Value = 15000
$ws.Cells.Item(121, 19).Value = 833

# Row 122
$ws.Cells.Item(122, 4).Value = 44390
$ws.Cells.Item(122, 13).Value = 140
$ws.Cells.Item(122, 14).Value = 12000
$ws.Cells.Item(122, 16).Value = 12571
$ws.Cells.Item(122, 19).Value = 698

# Row 123
$ws.Cells.Item(123, 4).Value = 44349
$ws.Cells.Item(123, 13).Value = 30
$ws.Cells.Item(123, 14).Value = 13000
$ws.Cells.Item(123, 16).Value = 13000
$ws.Cells.Item(123, 19).Value = 722

# New rows 124 and 125 (appended at the bottom, duplicating former rows 122-123 data).
# Column D needs the same date number-format as the rest of the 'Fecha' column.
# Row 124
$ws.Cells.Item(124, 1).Value = 10
$ws.Cells.Item(124, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(124, 3).Value = "La Araucanía"
$ws.Cells.Item(124, 4).Value = 44285
$ws.Cells.Item(124, 5).Value = 9
$ws.Cells.Item(124, 6).Value = "Fruta"
$ws.Cells.Item(124, 7).Value = 100104
$ws.Cells.Item(124, 8).Value = "Frutos de pepita"
$ws.Cells.Item(124, 9).Value = 100104003
$ws.Cells.Item(124, 10).Value = "Membrillo"
$ws.Cells.Item(124, 11).Value = "Champion"
$ws.Cells.Item(124, 12).Value = "Primera"
$ws.Cells.Item(124, 13).Value = 55
$ws.Cells.Item(124, 14).Value = 13000
$ws.Cells.Item(124, 15).Value = 13000
$ws.Cells.Item(124, 16).Value = 13000
$ws.Cells.Item(124, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(124, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(124, 19).Value = 722
$ws.Cells.Item(124, 20).Value = 18
$ws.Cells.Item(124, 4).NumberFormat = $ws.Cells.Item(29, 4).NumberFormat

# Row 125
$ws.Cells.Item(125, 1).Value = 10
$ws.Cells.Item(125, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(125, 3).Value = "La Araucanía"
$ws.Cells.Item(125, 4).Value = 44335
$ws.Cells.Item(125, 5).Value = 9
$ws.Cells.Item(125, 6).Value = "Fruta"
$ws.Cells.Item(125, 7).Value = 100104
$ws.Cells.Item(125, 8).Value = "Frutos de pepita"
$ws.Cells.Item(125, 9).Value = 100104003
$ws.Cells.Item(125, 10).Value = "Membrillo"
$ws.Cells.Item(125, 11).Value = "Champion"
$ws.Cells.Item(125, 12).Value = "Primera"
$ws.Cells.Item(125, 13).Value = 90
$ws.Cells.Item(125, 14).Value = 12000
$ws.Cells.Item(125, 15).Value = 13000
$ws.Cells.Item(125, 16).Value = 12556
$ws.Cells.Item(125, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(125, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(125, 19).Value = 698
$ws.Cells.Item(125, 20).Value = 18
$ws.Cells.Item(125, 4).NumberFormat = $ws.Cells.Item(29, 4).NumberFormat
